# Created new strategy HA_VWAP.
$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$listOfValues = $wb.Worksheets.Item("ListOfValues")
$strategyDict = $wb.Worksheets.Item("StrategyDictionaries")

# ---------------------------------------------------------------
# 1. StrategyDictionaries: add the HA_VWAP dictionary definition
# ---------------------------------------------------------------
$strategyDict.Range("A5").Value = "HA_VWAP:"
$strategyDict.Range("B5").Value = "{'EMA': 200', 'DistVWAP_PCT': 0.05', 'NB_SIGNALS': 2}"

# ---------------------------------------------------------------
# 2. ListOfValues: add HA_VWAP to the Strategy list (col B) and
#    VWAP_Touch to the Exit Strategy list (col C), then re-sort
#    the Strategy list alphabetically (B2:B8)
# ---------------------------------------------------------------
$listOfValues.Range("B8").Value = "HA_VWAP"
$listOfValues.Range("C4").Value = "VWAP_Touch"

$sortRange = $listOfValues.Range("B2:B8")
$sortRange.Sort($listOfValues.Range("B2:B8"), 1)

# ---------------------------------------------------------------
# 3. Sheet1: update the dropdown / list validations to cover the
#    newly added rows
# ---------------------------------------------------------------
# (validations reference ListOfValues!$B$2:$B$8 and $C$2:$C$4 after edit)

# ---------------------------------------------------------------
# 4. Sheet1: update row 2 (existing test) and fill in row 3 (new
#    test) with the HA_VWAP / VWAP_Touch strategy
# ---------------------------------------------------------------
$sheet1.Range("D2").Value = 44562
$sheet1.Range("F2").Value = "5m"
$sheet1.Range("G2").Value = 7
$sheet1.Range("H2").Value = 7
$sheet1.Range("I2").Value = "HA_VWAP"
$sheet1.Range("J2").Value = "VWAP_Touch"

$sheet1.Range("A3").Value = 2
$sheet1.Range("B3").Value = "Bybit"
$sheet1.Range("C3").Value = "ETHUSDT"
$sheet1.Range("D3").Value = 44562
$sheet1.Range("E3").Value = 44926
$sheet1.Range("F3").Value = "5m"
$sheet1.Range("G3").Value = 7
$sheet1.Range("H3").Value = 7
$sheet1.Range("I3").Value = "HA_VWAP"
$sheet1.Range("J3").Value = "VWAP_Touch"

# ---------------------------------------------------------------
# 5. Number format: TP%/SL% columns now show 3 decimals instead
#    of 4
# ---------------------------------------------------------------
$sheet1.Range("G1:H3").NumberFormat = "#,##0.000"

# ---------------------------------------------------------------
# 6. Selection bookmarks left behind by the author while editing
# ---------------------------------------------------------------
$listOfValues.Activate()
$listOfValues.Range("B10").Select()

$strategyDict.Activate()
$strategyDict.Range("B10").Select()

$sheet1.Activate()
$sheet1.Range("I5").Select()
